# Update "Forecast Comparison" sheet with a new Week_Start_Date column and
# corrected forecast output (per commit: "Update with Correct Forecast output").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date") - this shifts the existing
# ASIN / MyForecast / Amazon*Forecast / Product Title / is_holiday_week
# columns one to the right (B->C, C->D, ... I->J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels lose their leading zero (W01 -> W1, etc.) and each row gets
# its corresponding week-start date as plain text.
$weekStarts = @{
    2  = @{ Week = "W1";  Date = "2025-01-05" }
    3  = @{ Week = "W2";  Date = "2025-01-12" }
    4  = @{ Week = "W3";  Date = "2025-01-19" }
    5  = @{ Week = "W4";  Date = "2025-01-26" }
    6  = @{ Week = "W5";  Date = "2025-02-02" }
    7  = @{ Week = "W6";  Date = "2025-02-09" }
    8  = @{ Week = "W7";  Date = "2025-02-16" }
    9  = @{ Week = "W8";  Date = "2025-02-23" }
    10 = @{ Week = "W9";  Date = "2025-03-02" }
    11 = @{ Week = "W10"; Date = "2025-03-09" }
    12 = @{ Week = "W11"; Date = "2025-03-16" }
    13 = @{ Week = "W12"; Date = "2025-03-23" }
    14 = @{ Week = "W13"; Date = "2025-03-30" }
    15 = @{ Week = "W14"; Date = "2025-04-06" }
    16 = @{ Week = "W15"; Date = "2025-04-13" }
    17 = @{ Week = "W16"; Date = "2025-04-20" }
}

foreach ($row in $weekStarts.Keys) {
    $info = $weekStarts[$row]

    $ws.Range("A$row").Value = $info.Week

    # Force the date to be stored as plain text (matching the source data,
    # which keeps it as an inline string) instead of being auto-converted
    # to a date serial number.
    $dateCell = $ws.Range("B$row")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $info.Date
    $dateCell.ClearFormats()
}

# (Row 4's Amazon P70/P80/P90 Forecast figures shift naturally from the
# column insert above: old F4=1 -> new G4=1, old G4=1 -> new H4=1, and the
# newly vacated F4 inherits the old E4 value of 0.)

# is_holiday_week (now column J) becomes a proper boolean column.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}
